$d = $word.ActiveDocument
# Step 1: trim paragraph 6 text (up to bookmark), keep bookmark
$r1 = $d.Range(129, 238)
$r1.Delete()
$p6 = $d.Paragraphs.Item(6)
$r2 = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$r2.Delete()
Write-Output "After trim:"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "Para $i [$($p.Range.Start)-$($p.Range.End)] text=[$($p.Range.Text)]"
}
$p5 = $d.Paragraphs.Item(5)
$p6b = $d.Paragraphs.Item(6)
Write-Output "p5 end=$($p5.Range.End) p6 start=$($p6b.Range.Start) p6 end=$($p6b.Range.End)"

# try deleting just the p5 pilcrow char via a 2-element combined range (p5 mark + p6 mark), keep bookmark which sits between
$combined = $d.Range($p5.Range.End - 1, $p6b.Range.End)
Write-Output "combined [$($combined.Start)-$($combined.End)] text=[$($combined.Text)]"
$combined.Delete()
Write-Output "After combined delete:"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "Para $i [$($p.Range.Start)-$($p.Range.End)] text=[$($p.Range.Text)]"
}
